$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 16, shifting existing rows 16-77 down to 17-78.
$ws.Rows.Item(16).Insert()

# Populate the new row 16 with the new weekly price record.
$ws.Range("A16").Value = 10
$ws.Range("B16").Value = "Vega Modelo de Temuco"
$ws.Range("C16").Value = "La Araucanía"
$ws.Range("D16").Value = 45054
$ws.Range("E16").Value = 9
$ws.Range("F16").Value = 100112042
$ws.Range("G16").Value = "Locoto"
$ws.Range("H16").Value = "Sin especificar"
$ws.Range("I16").Value = "Primera"
$ws.Range("J16").Value = 90
$ws.Range("K16").Value = 4100
$ws.Range("L16").Value = 4100
$ws.Range("M16").Value = 4100
$ws.Range("N16").Value = "$/kilo"
$ws.Range("O16").Value = "Región de Arica y Parinacota"
$ws.Range("P16").Value = 4100
$ws.Range("Q16").Value = 1
$ws.Range("R16").Value = "Hortaliza"
